$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pt_max values in column E (rows 2-7) from 50 to 70
$ws.Range("E2:E7").Value = 70

# Update the selection to match the edited range
$ws.Range("E2:E7").Select()
